# Locate the paragraph that ends with "... journaler som ikke følger dette."
# i.e. the one whose text contains "ikke følger dette" and insert the new
# table + paragraphs right after it (before the following blank paragraph).

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*ikke følger dette*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find anchor paragraph."
}

$insertRange = $d.Range($target.Range.End, $target.Range.End)

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:tbl>
            <w:tblPr>
              <w:tblStyle w:val="Tabellrutenett"/>
              <w:tblW w:w="0" w:type="auto"/>
              <w:tblLayout w:type="fixed"/>
              <w:tblLook w:val="06A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="1" w:noVBand="1"/>
            </w:tblPr>
            <w:tblGrid>
              <w:gridCol w:w="4530"/>
              <w:gridCol w:w="4530"/>
            </w:tblGrid>
            <w:tr>
              <w:tc>
                <w:tcPr>
                  <w:tcW w:w="4530" w:type="dxa"/>
                </w:tcPr>
                <w:p>
                  <w:r>
                    <w:t>Journalstatus</w:t>
                  </w:r>
                </w:p>
              </w:tc>
              <w:tc>
                <w:tcPr>
                  <w:tcW w:w="4530" w:type="dxa"/>
                </w:tcPr>
                <w:p>
                  <w:r>
                    <w:t>Antall</w:t>
                  </w:r>
                </w:p>
              </w:tc>
            </w:tr>
          </w:tbl>
          <w:p/>
          <w:p>
            <w:r>
              <w:t>AND/OR</w:t>
            </w:r>
          </w:p>
          <w:p/>
          <w:p>
            <w:r>
              <w:rPr>
                <w:b/>
              </w:rPr>
              <w:t xml:space="preserve">Konsekvensvurdering: </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">Feilen anses liten, spesielt for et fysisk uttrekk, og godkjennes. </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($xmlFrag)

Write-Host "Inserted table and paragraphs after anchor paragraph."
